$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3.0
$ws.Range("F2").Value2 = 1.0
$ws.Range("G2").Value2 = 0.06437833333333333
$ws.Range("H2").Value2 = 0.193135
$ws.Range("I2").Value2 = 0.109187438766332
$ws.Range("J2").Value2 = 0.109187438766332
$ws.Range("M2").Value2 = 3.425446666666666
$ws.Range("N2").Value2 = 10.27634
$ws.Range("O2").Value2 = 0.6657953389778073
$ws.Range("P2").Value2 = 0.6657953389778073
$ws.Range("Q2").Value2 = 0.2205245473222222
$ws.Range("R2").Value2 = 1.9847209259
$ws.Range("S2").Value2 = 0.07269648780554862
$ws.Range("T2").Value2 = 0.07269648780554862

$ws.Range("E3").Value2 = 3.0
$ws.Range("F3").Value2 = 1.0
$ws.Range("G3").Value2 = 0.06437833333333333
$ws.Range("H3").Value2 = 0.193135
$ws.Range("I3").Value2 = 0.109187438766332
$ws.Range("J3").Value2 = 0.109187438766332
$ws.Range("O3").Value2 = 0.2094791321596951
$ws.Range("P3").Value2 = 0.2094791321596952
$ws.Range("Q3").Value2 = 0.06938361999333333
$ws.Range("R3").Value2 = 0.62445257994
$ws.Range("S3").Value2 = 0.02287248991551109
$ws.Range("T3").Value2 = 0.02287248991551109

$ws.Range("E4").Value2 = 3.0
$ws.Range("F4").Value2 = 1.0
$ws.Range("G4").Value2 = 0.06437833333333333
$ws.Range("H4").Value2 = 0.193135
$ws.Range("I4").Value2 = 0.109187438766332
$ws.Range("J4").Value2 = 0.109187438766332
$ws.Range("M4").Value2 = 0.62317
$ws.Range("N4").Value2 = 1.86951
$ws.Range("O4").Value2 = 0.1211239647746572
$ws.Range("P4").Value2 = 0.1211239647746572
$ws.Range("Q4").Value2 = 0.04011864598333333
$ws.Range("R4").Value2 = 0.36106781385
$ws.Range("S4").Value2 = 0.01322521548696824
$ws.Range("T4").Value2 = 0.01322521548696824

$ws.Range("E5").Value2 = 3.0
$ws.Range("F5").Value2 = 1.0
$ws.Range("G5").Value2 = 0.06437833333333333
$ws.Range("H5").Value2 = 0.193135
$ws.Range("I5").Value2 = 0.109187438766332
$ws.Range("J5").Value2 = 0.109187438766332
$ws.Range("M5").Value2 = 0.01852966666666667
$ws.Range("N5").Value2 = 0.055589
$ws.Range("O5").Value2 = 0.003601564087840353
$ws.Range("P5").Value2 = 0.003601564087840353
$ws.Range("Q5").Value2 = 0.001192909057222222
$ws.Range("R5").Value2 = 0.010736181515
$ws.Range("S5").Value2 = 0.000393245558304089
$ws.Range("T5").Value2 = 0.000393245558304089

$ws.Range("G6").Value2 = 0.4788196666666666
$ws.Range("I6").Value2 = 0.8120914339857952
$ws.Range("J6").Value2 = 0.8120914339857951
$ws.Range("M6").Value2 = 3.425446666666666
$ws.Range("N6").Value2 = 10.27634
$ws.Range("O6").Value2 = 0.6657953389778073
$ws.Range("P6").Value2 = 0.6657953389778073
$ws.Range("Q6").Value2 = 1.640171231117778
$ws.Range("R6").Value2 = 14.76154108006
$ws.Range("S6").Value2 = 0.5406866915715461
$ws.Range("T6").Value2 = 0.540686691571546

$ws.Range("G7").Value2 = 0.4788196666666666
$ws.Range("I7").Value2 = 0.8120914339857952
$ws.Range("J7").Value2 = 0.8120914339857951
$ws.Range("O7").Value2 = 0.2094791321596951
$ws.Range("P7").Value2 = 0.2094791321596952
$ws.Range("Q7").Value2 = 0.5160469381106666
$ws.Range("R7").Value2 = 4.644422442995999
$ws.Range("S7").Value2 = 0.1701162088256667
$ws.Range("T7").Value2 = 0.1701162088256667

$ws.Range("G8").Value2 = 0.4788196666666666
$ws.Range("I8").Value2 = 0.8120914339857952
$ws.Range("J8").Value2 = 0.8120914339857951
$ws.Range("M8").Value2 = 0.62317
$ws.Range("N8").Value2 = 1.86951
$ws.Range("O8").Value2 = 0.1211239647746572
$ws.Range("P8").Value2 = 0.1211239647746572
$ws.Range("Q8").Value2 = 0.2983860516766667
$ws.Range("R8").Value2 = 2.68547446509
$ws.Range("S8").Value2 = 0.0983637342438963
$ws.Range("T8").Value2 = 0.09836373424389629

$ws.Range("G9").Value2 = 0.4788196666666666
$ws.Range("I9").Value2 = 0.8120914339857952
$ws.Range("J9").Value2 = 0.8120914339857951
$ws.Range("M9").Value2 = 0.01852966666666667
$ws.Range("N9").Value2 = 0.055589
$ws.Range("O9").Value2 = 0.003601564087840353
$ws.Range("P9").Value2 = 0.003601564087840353
$ws.Range("Q9").Value2 = 0.008872368816777778
$ws.Range("R9").Value2 = 0.079851319351
$ws.Range("S9").Value2 = 0.002924799344686015
$ws.Range("T9").Value2 = 0.002924799344686014

$ws.Range("G10").Value2 = 0.042481
$ws.Range("H10").Value2 = 0.127443
$ws.Range("I10").Value2 = 0.07204895414449818
$ws.Range("J10").Value2 = 0.07204895414449818
$ws.Range("M10").Value2 = 3.425446666666666
$ws.Range("N10").Value2 = 10.27634
$ws.Range("O10").Value2 = 0.6657953389778073
$ws.Range("P10").Value2 = 0.6657953389778073
$ws.Range("Q10").Value2 = 0.1455163998466666
$ws.Range("R10").Value2 = 1.30964759862
$ws.Range("S10").Value2 = 0.04796985784763266
$ws.Range("T10").Value2 = 0.04796985784763266

$ws.Range("G11").Value2 = 0.042481
$ws.Range("H11").Value2 = 0.127443
$ws.Range("I11").Value2 = 0.07204895414449818
$ws.Range("J11").Value2 = 0.07204895414449818
$ws.Range("O11").Value2 = 0.2094791321596951
$ws.Range("P11").Value2 = 0.2094791321596952
$ws.Range("Q11").Value2 = 0.04578381278799999
$ws.Range("R11").Value2 = 0.412054315092
$ws.Range("S11").Value2 = 0.01509275238720315
$ws.Range("T11").Value2 = 0.01509275238720315

$ws.Range("G12").Value2 = 0.042481
$ws.Range("H12").Value2 = 0.127443
$ws.Range("I12").Value2 = 0.07204895414449818
$ws.Range("J12").Value2 = 0.07204895414449818
$ws.Range("M12").Value2 = 0.62317
$ws.Range("N12").Value2 = 1.86951
$ws.Range("O12").Value2 = 0.1211239647746572
$ws.Range("P12").Value2 = 0.1211239647746572
$ws.Range("Q12").Value2 = 0.02647288477
$ws.Range("R12").Value2 = 0.23825596293
$ws.Range("S12").Value2 = 0.008726854983849088
$ws.Range("T12").Value2 = 0.008726854983849088

$ws.Range("G13").Value2 = 0.042481
$ws.Range("H13").Value2 = 0.127443
$ws.Range("I13").Value2 = 0.07204895414449818
$ws.Range("J13").Value2 = 0.07204895414449818
$ws.Range("M13").Value2 = 0.01852966666666667
$ws.Range("N13").Value2 = 0.055589
$ws.Range("O13").Value2 = 0.003601564087840353
$ws.Range("P13").Value2 = 0.003601564087840353
$ws.Range("Q13").Value2 = 0.0007871587696666666
$ws.Range("R13").Value2 = 0.007084428927
$ws.Range("S13").Value2 = 0.000259488925813281
$ws.Range("T13").Value2 = 0.000259488925813281

$ws.Range("E14").Value2 = 1.0
$ws.Range("F14").Value2 = 0.3333333333333333
$ws.Range("G14").Value2 = 0.003934
$ws.Range("H14").Value2 = 0.011802
$ws.Range("I14").Value2 = 0.006672173103374587
$ws.Range("J14").Value2 = 0.006672173103374586
$ws.Range("M14").Value2 = 3.425446666666666
$ws.Range("N14").Value2 = 10.27634
$ws.Range("O14").Value2 = 0.6657953389778073
$ws.Range("P14").Value2 = 0.6657953389778073
$ws.Range("Q14").Value2 = 0.01347570718666666
$ws.Range("R14").Value2 = 0.12128136468
$ws.Range("S14").Value2 = 0.004442301753079891
$ws.Range("T14").Value2 = 0.004442301753079891

$ws.Range("E15").Value2 = 1.0
$ws.Range("F15").Value2 = 0.3333333333333333
$ws.Range("G15").Value2 = 0.003934
$ws.Range("H15").Value2 = 0.011802
$ws.Range("I15").Value2 = 0.006672173103374587
$ws.Range("J15").Value2 = 0.006672173103374586
$ws.Range("O15").Value2 = 0.2094791321596951
$ws.Range("P15").Value2 = 0.2094791321596952
$ws.Range("Q15").Value2 = 0.004239860632
$ws.Range("R15").Value2 = 0.038158745688
$ws.Range("S15").Value2 = 0.001397681031314168
$ws.Range("T15").Value2 = 0.001397681031314168

$ws.Range("E16").Value2 = 1.0
$ws.Range("F16").Value2 = 0.3333333333333333
$ws.Range("G16").Value2 = 0.003934
$ws.Range("H16").Value2 = 0.011802
$ws.Range("I16").Value2 = 0.006672173103374587
$ws.Range("J16").Value2 = 0.006672173103374586
$ws.Range("M16").Value2 = 0.62317
$ws.Range("N16").Value2 = 1.86951
$ws.Range("O16").Value2 = 0.1211239647746572
$ws.Range("P16").Value2 = 0.1211239647746572
$ws.Range("Q16").Value2 = 0.00245155078
$ws.Range("R16").Value2 = 0.02206395702
$ws.Range("S16").Value2 = 0.0008081600599435585
$ws.Range("T16").Value2 = 0.0008081600599435584

$ws.Range("E17").Value2 = 1.0
$ws.Range("F17").Value2 = 0.3333333333333333
$ws.Range("G17").Value2 = 0.003934
$ws.Range("H17").Value2 = 0.011802
$ws.Range("I17").Value2 = 0.006672173103374587
$ws.Range("J17").Value2 = 0.006672173103374586
$ws.Range("M17").Value2 = 0.01852966666666667
$ws.Range("N17").Value2 = 0.055589
$ws.Range("O17").Value2 = 0.003601564087840353
$ws.Range("P17").Value2 = 0.003601564087840353
$ws.Range("Q17").Value2 = 0.00007289570866666666
$ws.Range("R17").Value2 = 0.000656061378
$ws.Range("S17").Value2 = 0.00002403025903696823
$ws.Range("T17").Value2 = 0.00002403025903696822
